$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 57.319636
$ws.Range("H2").Value = 171.958908
$ws.Range("I2").Value = 0.5476981520382651
$ws.Range("J2").Value = 0.5476981520382651
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 3.438907666666667
$ws.Range("N2").Value = 10.316723
$ws.Range("O2").Value = 0.05825422340060618
$ws.Range("P2").Value = 0.05825422340060618
$ws.Range("Q2").Value = 197.1169356909427
$ws.Range("R2").Value = 1774.052421218484
$ws.Range("S2").Value = 0.03190573050493627
$ws.Range("T2").Value = 0.03190573050493627

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 57.319636
$ws.Range("H3").Value = 171.958908
$ws.Range("I3").Value = 0.5476981520382651
$ws.Range("J3").Value = 0.5476981520382651
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 10.383857
$ws.Range("N3").Value = 31.151571
$ws.Range("O3").Value = 0.1758999031294962
$ws.Range("P3").Value = 0.1758999031294962
$ws.Range("Q3").Value = 595.198903516052
$ws.Range("R3").Value = 5356.790131644469
$ws.Range("S3").Value = 0.0963400518877349
$ws.Range("T3").Value = 0.0963400518877349

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 57.319636
$ws.Range("H4").Value = 171.958908
$ws.Range("I4").Value = 0.5476981520382651
$ws.Range("J4").Value = 0.5476981520382651
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 25.34077833333333
$ws.Range("N4").Value = 76.022335
$ws.Range("O4").Value = 0.4292663558501786
$ws.Range("P4").Value = 0.4292663558501786
$ws.Range("Q4").Value = 1452.524190023353
$ws.Range("R4").Value = 13072.71771021018
$ws.Range("S4").Value = 0.2351083898313431
$ws.Range("T4").Value = 0.2351083898313432

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 57.319636
$ws.Range("H5").Value = 171.958908
$ws.Range("I5").Value = 0.5476981520382651
$ws.Range("J5").Value = 0.5476981520382651
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 19.86921833333334
$ws.Range("N5").Value = 59.60765500000001
$ws.Range("O5").Value = 0.336579517619719
$ws.Range("P5").Value = 0.336579517619719
$ws.Range("Q5").Value = 1138.896362471193
$ws.Range("R5").Value = 10250.06726224074
$ws.Range("S5").Value = 0.1843439798142508
$ws.Range("T5").Value = 0.1843439798142508

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 5.975184333333334
$ws.Range("H6").Value = 17.925553
$ws.Range("I6").Value = 0.05709382762749331
$ws.Range("J6").Value = 0.05709382762749331
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 3.438907666666667
$ws.Range("N6").Value = 10.316723
$ws.Range("O6").Value = 0.05825422340060618
$ws.Range("P6").Value = 0.05825422340060618
$ws.Range("Q6").Value = 20.54810721364656
$ws.Range("R6").Value = 184.932964922819
$ws.Range("S6").Value = 0.003325956589407696
$ws.Range("T6").Value = 0.003325956589407696

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 5.975184333333334
$ws.Range("H7").Value = 17.925553
$ws.Range("I7").Value = 0.05709382762749331
$ws.Range("J7").Value = 0.05709382762749331
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 10.383857
$ws.Range("N7").Value = 31.151571
$ws.Range("O7").Value = 0.1758999031294962
$ws.Range("P7").Value = 0.1758999031294962
$ws.Range("Q7").Value = 62.04545966597367
$ws.Range("R7").Value = 558.409136993763
$ws.Range("S7").Value = 0.01004279874896822
$ws.Range("T7").Value = 0.01004279874896822

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 5.975184333333334
$ws.Range("H8").Value = 17.925553
$ws.Range("I8").Value = 0.05709382762749331
$ws.Range("J8").Value = 0.05709382762749331
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 25.34077833333333
$ws.Range("N8").Value = 76.022335
$ws.Range("O8").Value = 0.4292663558501786
$ws.Range("P8").Value = 0.4292663558501786
$ws.Range("Q8").Value = 151.4158216918061
$ws.Range("R8").Value = 1362.742395226255
$ws.Range("S8").Value = 0.0245084593271923
$ws.Range("T8").Value = 0.0245084593271923

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 5.975184333333334
$ws.Range("H9").Value = 17.925553
$ws.Range("I9").Value = 0.05709382762749331
$ws.Range("J9").Value = 0.05709382762749331
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 19.86921833333334
$ws.Range("N9").Value = 59.60765500000001
$ws.Range("O9").Value = 0.336579517619719
$ws.Range("P9").Value = 0.336579517619719
$ws.Range("Q9").Value = 118.7222421009128
$ws.Range("R9").Value = 1068.500178908215
$ws.Range("S9").Value = 0.01921661296192508
$ws.Range("T9").Value = 0.01921661296192509

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 32.32302533333333
$ws.Range("H10").Value = 96.969076
$ws.Range("I10").Value = 0.3088515991858827
$ws.Range("J10").Value = 0.3088515991858827
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 3.438907666666667
$ws.Range("N10").Value = 10.316723
$ws.Range("O10").Value = 0.05825422340060618
$ws.Range("P10").Value = 0.05825422340060618
$ws.Range("Q10").Value = 111.1558996286609
$ws.Range("R10").Value = 1000.403096657948
$ws.Range("S10").Value = 0.01799191005660889
$ws.Range("T10").Value = 0.01799191005660889

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 32.32302533333333
$ws.Range("H11").Value = 96.969076
$ws.Range("I11").Value = 0.3088515991858827
$ws.Range("J11").Value = 0.3088515991858827
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 10.383857
$ws.Range("N11").Value = 31.151571
$ws.Range("O11").Value = 0.1758999031294962
$ws.Range("P11").Value = 0.1758999031294962
$ws.Range("Q11").Value = 335.6376728687107
$ws.Range("R11").Value = 3020.739055818396
$ws.Range("S11").Value = 0.05432696637818675
$ws.Range("T11").Value = 0.05432696637818675

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 32.32302533333333
$ws.Range("H12").Value = 96.969076
$ws.Range("I12").Value = 0.3088515991858827
$ws.Range("J12").Value = 0.3088515991858827
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 25.34077833333333
$ws.Range("N12").Value = 76.022335
$ws.Range("O12").Value = 0.4292663558501786
$ws.Range("P12").Value = 0.4292663558501786
$ws.Range("Q12").Value = 819.0906200347177
$ws.Range("R12").Value = 7371.81558031246
$ws.Range("S12").Value = 0.1325796004810239
$ws.Range("T12").Value = 0.1325796004810239

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 32.32302533333333
$ws.Range("H13").Value = 96.969076
$ws.Range("I13").Value = 0.3088515991858827
$ws.Range("J13").Value = 0.3088515991858827
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 19.86921833333334
$ws.Range("N13").Value = 59.60765500000001
$ws.Range("O13").Value = 0.336579517619719
$ws.Range("P13").Value = 0.336579517619719
$ws.Range("Q13").Value = 642.2332475418646
$ws.Range("R13").Value = 5780.099227876781
$ws.Range("S13").Value = 0.1039531222700632
$ws.Range("T13").Value = 0.1039531222700632

# Row 14
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 9.037676333333334
$ws.Range("H14").Value = 27.113029
$ws.Range("I14").Value = 0.08635642114835883
$ws.Range("J14").Value = 0.08635642114835884
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 3.438907666666667
$ws.Range("N14").Value = 10.316723
$ws.Range("O14").Value = 0.05825422340060618
$ws.Range("P14").Value = 0.05825422340060618
$ws.Range("Q14").Value = 31.07973443155189
$ws.Range("R14").Value = 279.717609883967
$ws.Range("S14").Value = 0.005030626249653327
$ws.Range("T14").Value = 0.005030626249653328

# Row 15
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 9.037676333333334
$ws.Range("H15").Value = 27.113029
$ws.Range("I15").Value = 0.08635642114835883
$ws.Range("J15").Value = 0.08635642114835884
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 10.383857
$ws.Range("N15").Value = 31.151571
$ws.Range("O15").Value = 0.1758999031294962
$ws.Range("P15").Value = 0.1758999031294962
$ws.Range("Q15").Value = 93.84593865761768
$ws.Range("R15").Value = 844.613447918559
$ws.Range("S15").Value = 0.01519008611460629
$ws.Range("T15").Value = 0.01519008611460629

# Row 16
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 9.037676333333334
$ws.Range("H16").Value = 27.113029
$ws.Range("I16").Value = 0.08635642114835883
$ws.Range("J16").Value = 0.08635642114835884
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 25.34077833333333
$ws.Range("N16").Value = 76.022335
$ws.Range("O16").Value = 0.4292663558501786
$ws.Range("P16").Value = 0.4292663558501786
$ws.Range("Q16").Value = 229.0217526114128
$ws.Range("R16").Value = 2061.195773502715
$ws.Range("S16").Value = 0.03706990621061929
$ws.Range("T16").Value = 0.0370699062106193

# Row 17
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 9.037676333333334
$ws.Range("H17").Value = 27.113029
$ws.Range("I17").Value = 0.08635642114835883
$ws.Range("J17").Value = 0.08635642114835884
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 19.86921833333334
$ws.Range("N17").Value = 59.60765500000001
$ws.Range("O17").Value = 0.336579517619719
$ws.Range("P17").Value = 0.336579517619719
$ws.Range("Q17").Value = 179.5715642929995
$ws.Range("R17").Value = 1616.144078636995
$ws.Range("S17").Value = 0.02906580257347991
$ws.Range("T17").Value = 0.02906580257347992
